$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AG2").Value = 0
$ws.Range("AJ2").Value = 12731182
$ws.Range("D2:AF2").ClearContents()
$ws.Range("AH2:AI2").ClearContents()

# Row 3
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AJ3").Value = 17089768
$ws.Range("D3:AF3").ClearContents()
$ws.Range("AI3").ClearContents()

# Row 4
$ws.Range("D4").Value = 1556
$ws.Range("E4").Value = 449
$ws.Range("F4").Value = 449
$ws.Range("G4").Value = 534
$ws.Range("H4").Value = 487
$ws.Range("I4").Value = 488
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 4370
$ws.Range("L4").Value = 587
$ws.Range("M4").Value = 3782
$ws.Range("N4").Value = 3782
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 86
$ws.Range("Q4").Value = 543
$ws.Range("R4").Value = -1059
$ws.Range("S4").Value = -64
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 538
$ws.Range("V4").Value = 435
$ws.Range("W4").Value = 28.82
$ws.Range("X4").Value = 31.3
$ws.Range("AA4").Value = 15.53
$ws.Range("AB4").Value = 4612.25
$ws.Range("AC4").Value = 2847
$ws.Range("AD4").Value = 14
$ws.Range("AE4").Value = 23957
$ws.Range("AF4").Value = 1.66
$ws.Range("AG4").Value = 650
$ws.Range("AH4").Value = 1.63
$ws.Range("AI4").Value = 21.07
$ws.Range("AJ4").Value = 17184556
$ws.Range("Y4:Z4").ClearContents()

# Row 5
$ws.Range("D5").Value = 3193
$ws.Range("E5").Value = 810
$ws.Range("F5").Value = 810
$ws.Range("G5").Value = 383
$ws.Range("H5").Value = 356
$ws.Range("I5").Value = 354
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 9972
$ws.Range("L5").Value = 5926
$ws.Range("M5").Value = 4046
$ws.Range("N5").Value = 4046
$ws.Range("P5").Value = 88
$ws.Range("Q5").Value = 824
$ws.Range("R5").Value = -5968
$ws.Range("S5").Value = 5246
$ws.Range("T5").Value = 6
$ws.Range("U5").Value = 818
$ws.Range("V5").Value = 5492
$ws.Range("W5").Value = 25.37
$ws.Range("X5").Value = 11.16
$ws.Range("Y5").Value = 9.039999999999999
$ws.Range("Z5").Value = 4.97
$ws.Range("AA5").Value = 146.45
$ws.Range("AB5").Value = 5192.3
$ws.Range("AC5").Value = 2027
$ws.Range("AD5").Value = 25.5
$ws.Range("AE5").Value = 24982
$ws.Range("AF5").Value = 2.07
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 0.68
$ws.Range("AI5").Value = 16.03
$ws.Range("AJ5").Value = 17595514
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 4830
$ws.Range("E6").Value = 1358
$ws.Range("F6").Value = 1358
$ws.Range("G6").Value = 1174
$ws.Range("H6").Value = 878
$ws.Range("I6").Value = 878
$ws.Range("K6").Value = 9968
$ws.Range("L6").Value = 4785
$ws.Range("M6").Value = 5183
$ws.Range("N6").Value = 5183
$ws.Range("P6").Value = 91
$ws.Range("Q6").Value = 1563
$ws.Range("R6").Value = 34
$ws.Range("S6").Value = -1431
$ws.Range("T6").Value = 7
$ws.Range("U6").Value = 1556
$ws.Range("V6").Value = 4238
$ws.Range("W6").Value = 28.11
$ws.Range("X6").Value = 18.17
$ws.Range("Y6").Value = 19.02
$ws.Range("Z6").Value = 8.800000000000001
$ws.Range("AA6").Value = 92.33
$ws.Range("AB6").Value = 6059.62
$ws.Range("AC6").Value = 4869
$ws.Range("AD6").Value = 12.28
$ws.Range("AE6").Value = 30760
$ws.Range("AF6").Value = 1.94
$ws.Range("AI6").Value = 6.72
$ws.Range("AJ6").Value = 18247298
$ws.Range("AG6:AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 5197
$ws.Range("E7").Value = 1530
$ws.Range("G7").Value = 1460
$ws.Range("H7").Value = 1155
$ws.Range("I7").Value = 1155
$ws.Range("K7").Value = 10682
$ws.Range("L7").Value = 4247
$ws.Range("M7").Value = 6434
$ws.Range("N7").Value = 6476
$ws.Range("P7").Value = 91
$ws.Range("Q7").Value = 1532
$ws.Range("R7").Value = -345
$ws.Range("S7").Value = -969
$ws.Range("T7").Value = 16
$ws.Range("U7").Value = 1488
$ws.Range("W7").Value = 29.45
$ws.Range("X7").Value = 22.23
$ws.Range("Y7").Value = 19.82
$ws.Range("Z7").Value = 11.19
$ws.Range("AA7").Value = 66.01000000000001
$ws.Range("AC7").Value = 6327
$ws.Range("AD7").Value = 7.47
$ws.Range("AE7").Value = 38405
$ws.Range("AF7").Value = 1.23
$ws.Range("AG7").Value = 343
$ws.Range("AH7").Value = 0.73
$ws.Range("AI7").Value = 5.43

# Row 8
$ws.Range("D8").Value = 5546
$ws.Range("E8").Value = 1646
$ws.Range("G8").Value = 1607
$ws.Range("H8").Value = 1258
$ws.Range("I8").Value = 1258
$ws.Range("K8").Value = 12070
$ws.Range("L8").Value = 4312
$ws.Range("M8").Value = 7757
$ws.Range("N8").Value = 7846
$ws.Range("P8").Value = 91
$ws.Range("Q8").Value = 1728
$ws.Range("R8").Value = -298
$ws.Range("S8").Value = -124
$ws.Range("T8").Value = 23
$ws.Range("U8").Value = 1645
$ws.Range("W8").Value = 29.68
$ws.Range("X8").Value = 22.69
$ws.Range("Y8").Value = 17.57
$ws.Range("Z8").Value = 11.06
$ws.Range("AA8").Value = 55.59
$ws.Range("AC8").Value = 6891
$ws.Range("AD8").Value = 6.86
$ws.Range("AE8").Value = 46532
$ws.Range("AF8").Value = 1.02
$ws.Range("AG8").Value = 343
$ws.Range("AH8").Value = 0.73
$ws.Range("AI8").Value = 4.98

# Row 9
$ws.Range("D9").Value = 5959
$ws.Range("E9").Value = 1819
$ws.Range("G9").Value = 1835
$ws.Range("H9").Value = 1462
$ws.Range("I9").Value = 1462
$ws.Range("K9").Value = 13736
$ws.Range("L9").Value = 4321
$ws.Range("M9").Value = 9415
$ws.Range("N9").Value = 9415
$ws.Range("P9").Value = 91
$ws.Range("Q9").Value = 1904
$ws.Range("R9").Value = -401
$ws.Range("S9").Value = -145
$ws.Range("T9").Value = 30
$ws.Range("U9").Value = 1813
$ws.Range("W9").Value = 30.53
$ws.Range("X9").Value = 24.54
$ws.Range("Y9").Value = 16.94
$ws.Range("Z9").Value = 11.33
$ws.Range("AA9").Value = 45.9
$ws.Range("AC9").Value = 8008
$ws.Range("AD9").Value = 5.9
$ws.Range("AE9").Value = 55835
$ws.Range("AF9").Value = 0.85
$ws.Range("AG9").Value = 341
$ws.Range("AH9").Value = 0.72
$ws.Range("AI9").Value = 4.26
